# Apply the "drivers" sheet rework described by the diff:
#  - rename sheet "Conducteurs" -> "drivers"
#  - add a new header "vehicle" in G1
#  - clear the vehicle.matricule / vehicle.model values on rows 5 and 6 (E5:F5, E6:F6)
#  - append 7 new rows (7-13) of placeholder driver data in columns A:D
#  - dimension / ignoredErrors ranges follow automatically from the used range

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the worksheet.
$ws.Name = "drivers"

# New header column G.
$ws.Range("G1").Value = "vehicle"

# Remove the vehicle matricule/model data for rows 5 and 6.
$ws.Range("E5").ClearContents()
$ws.Range("F5").ClearContents()
$ws.Range("E6").ClearContents()
$ws.Range("F6").ClearContents()

# Append new placeholder rows 7..13 (ids 6..12).
$newRows = @(
    @(6,  "sdsad", "sdasd", "23213213"),
    @(7,  "sdsad", "sdasd", "23213213"),
    @(8,  "sdsad", "sdasd", "23213213"),
    @(9,  "sdsad", "sdasd", "23213213"),
    @(10, "sdsad", "sdasd", "23213213"),
    @(11, "sdsad", "sdasd", "23213213"),
    @(12, "sdsad", "sdasd", "23213213")
)

# Column D holds digit-only values that must stay text (matches the
# existing "phone"-style text-as-number cells in column D above), so
# format it as text before assigning.
$ws.Range("D7:D13").NumberFormat = "@"

$rowIndex = 7
foreach ($row in $newRows) {
    $ws.Cells.Item($rowIndex, 1).Value = $row[0]
    $ws.Cells.Item($rowIndex, 2).Value = $row[1]
    $ws.Cells.Item($rowIndex, 3).Value = $row[2]
    $ws.Cells.Item($rowIndex, 4).Value = $row[3]
    $rowIndex++
}
